$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.408.47"
$ws.Range("E2").Value = "  +0.10%  "

$ws.Range("D3").Value = "1.849.52"
$ws.Range("E3").Value = "  +0.17%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.59"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.29%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6280"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.23%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07628"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.36%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2910"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.66%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.69"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.84%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07734"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.07%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.031"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.65%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6790"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.27%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.00001059"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -3.13%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "83.20"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.50%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.157"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.13%  "

$ws.Range("D17").Value = "29.454.16"
$ws.Range("E17").Value = "  +0.18%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "226.54"
$ws.Range("D18").ClearFormats()

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.31"
$ws.Range("D19").ClearFormats()

$ws.Range("E20").Value = "  +0.02%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.498"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.11%  "

$ws.Range("E22").Value = "  +0.04%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "157.93"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.54%  "

$ws.Range("E24").Value = "  -0.86%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.397"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.20%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "17.68"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.36%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.395"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +6.56%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.462"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.15%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.05596"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.24%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.123"
$ws.Range("D30").ClearFormats()

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.064"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.76%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.836"
$ws.Range("D32").ClearFormats()

$ws.Range("E33").Value = "  +0.51%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.6956"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.96%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.590"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.25%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.01802"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.25%  "

$ws.Range("D37").Value = "1.229.41"
$ws.Range("E37").Value = "  -0.14%  "

$ws.Range("E38").Value = "  -1.67%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.389"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.89%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9025"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.11%  "

$ws.Range("E41").Value = "  +0.13%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "101.55"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.01%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "65.91"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.14%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.163"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.22%  "

$ws.Range("E45").Value = "  -3.29%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4009"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.10%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.008"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.70%  "

$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.678"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.05%  "

$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1144"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.04%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05703"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.10%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4632"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.17%  "
